$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "30.147.32"
Set-TextValue $ws.Range("E2") "  -0.61%  "
Set-TextValue $ws.Range("D3") "1.913.40"
Set-TextValue $ws.Range("E3") "  -1.11%  "
Set-TextValue $ws.Range("D4") "1.001"
Set-TextValue $ws.Range("E4") "  +0.08%  "
Set-TextValue $ws.Range("D5") "0.7385"
Set-TextValue $ws.Range("E5") "  -2.62%  "
Set-TextValue $ws.Range("D6") "243.84"
Set-TextValue $ws.Range("E6") "  -0.43%  "
Set-TextValue $ws.Range("E7") "  +0.09%  "
Set-TextValue $ws.Range("D8") "0.3125"
Set-TextValue $ws.Range("E8") "  -1.77%  "
Set-TextValue $ws.Range("D9") "26.72"
Set-TextValue $ws.Range("E9") "  -3.41%  "
Set-TextValue $ws.Range("D10") "0.06972"
Set-TextValue $ws.Range("E10") "  -0.32%  "
Set-TextValue $ws.Range("D11") "0.7802"
Set-TextValue $ws.Range("E11") "  +0.14%  "
Set-TextValue $ws.Range("D12") "0.07984"
Set-TextValue $ws.Range("E12") "  -0.24%  "
Set-TextValue $ws.Range("D13") "1.907.44"
Set-TextValue $ws.Range("E13") "  -1.44%  "
Set-TextValue $ws.Range("D14") "5.290"
Set-TextValue $ws.Range("E14") "  -1.12%  "
Set-TextValue $ws.Range("D15") "92.29"
Set-TextValue $ws.Range("E15") "  -2.15%  "
Set-TextValue $ws.Range("D16") "14.38"
Set-TextValue $ws.Range("E16") "  -0.17%  "
Set-TextValue $ws.Range("D17") "30.164.29"
Set-TextValue $ws.Range("E17") "  -0.51%  "
Set-TextValue $ws.Range("D18") "5.920"
Set-TextValue $ws.Range("E18") "  +2.93%  "
Set-TextValue $ws.Range("D19") "242.44"
Set-TextValue $ws.Range("E19") "  -4.12%  "
Set-TextValue $ws.Range("D20") "0.000007830"
Set-TextValue $ws.Range("E20") "  -1.20%  "
Set-TextValue $ws.Range("D21") "1.000"
Set-TextValue $ws.Range("E21") "  +0.11%  "
Set-TextValue $ws.Range("D22") "2.139.98"
Set-TextValue $ws.Range("E22") "  -2.36%  "
Set-TextValue $ws.Range("D23") "1.001"
Set-TextValue $ws.Range("E23") "  +0.07%  "
Set-TextValue $ws.Range("D24") "7.172"
Set-TextValue $ws.Range("E24") "  +7.62%  "
Set-TextValue $ws.Range("D25") "9.428"
Set-TextValue $ws.Range("E25") "  -0.40%  "
Set-TextValue $ws.Range("D26") "168.71"
Set-TextValue $ws.Range("E26") "  +1.92%  "
Set-TextValue $ws.Range("D27") "19.12"
Set-TextValue $ws.Range("E27") "  +0.70%  "
Set-TextValue $ws.Range("D28") "0.1283"
Set-TextValue $ws.Range("E28") "  -3.88%  "
Set-TextValue $ws.Range("D29") "2.071"
Set-TextValue $ws.Range("E29") "  -5.71%  "
Set-TextValue $ws.Range("D30") "1.355"
Set-TextValue $ws.Range("E30") "  -0.67%  "
Set-TextValue $ws.Range("D31") "1.547"
Set-TextValue $ws.Range("E31") "  +2.00%  "
Set-TextValue $ws.Range("D32") "4.347"
Set-TextValue $ws.Range("E32") "  -1.15%  "
Set-TextValue $ws.Range("D33") "4.105"
Set-TextValue $ws.Range("E33") "  -0.41%  "
Set-TextValue $ws.Range("D34") "0.05169"
Set-TextValue $ws.Range("D35") "1.300"
Set-TextValue $ws.Range("E35") "  +1.25%  "
Set-TextValue $ws.Range("D36") "0.7502"
Set-TextValue $ws.Range("E36") "  -0.07%  "
Set-TextValue $ws.Range("D37") "2.733"
Set-TextValue $ws.Range("E37") "  -1.35%  "
Set-TextValue $ws.Range("D38") "0.01947"
Set-TextValue $ws.Range("E38") "  -0.56%  "
Set-TextValue $ws.Range("D39") "2.800"
Set-TextValue $ws.Range("E39") "  +0.04%  "
Set-TextValue $ws.Range("D40") "6.376"
Set-TextValue $ws.Range("E40") "  -0.49%  "
Set-TextValue $ws.Range("D41") "75.18"
Set-TextValue $ws.Range("E41") "  -2.93%  "
Set-TextValue $ws.Range("D42") "0.4513"
Set-TextValue $ws.Range("E42") "  +1.28%  "
Set-TextValue $ws.Range("D43") "1.963"
Set-TextValue $ws.Range("E43") "  -0.08%  "
Set-TextValue $ws.Range("D44") "7.874"
Set-TextValue $ws.Range("E44") "  +5.34%  "
Set-TextValue $ws.Range("D45") "1.001"
Set-TextValue $ws.Range("E45") "  +0.15%  "
Set-TextValue $ws.Range("D46") "0.8389"
Set-TextValue $ws.Range("E46") "  +0.65%  "
Set-TextValue $ws.Range("D47") "9.948"
Set-TextValue $ws.Range("E47") "  +1.87%  "
Set-TextValue $ws.Range("D48") "101.50"
Set-TextValue $ws.Range("E48") "  +0.74%  "
Set-TextValue $ws.Range("D49") "37.28"
Set-TextValue $ws.Range("E49") "  -0.40%  "
Set-TextValue $ws.Range("D50") "2.048.44"
Set-TextValue $ws.Range("E50") "  -1.62%  "
Set-TextValue $ws.Range("D51") "0.1196"
Set-TextValue $ws.Range("E51") "  +1.69%  "
